# Rewrites the three "OBJETIVOS ESPECIFICOS" bullet paragraphs
# (Implementar / Organizar / Facilitar) into plain, non-numbered
# paragraphs with an inline "• " marker, Segoe UI font, light-grey
# shading, expanded body copy, and drops the final empty paragraph.

$d = $word.ActiveDocument

# Locate the first paragraph of the bulleted block ("Implementar...")
# and treat everything from there through the end of the document
# (which also swallows the trailing empty paragraph) as the block to
# replace.
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Implementar*") {
        $startPara = $p
        break
    }
}

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$target = $d.Range($startPara.Range.Start, $lastPara.Range.End)

$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:shd w:val="clear" w:color="auto" w:fill="F7F7F8"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:shd w:val="clear" w:color="auto" w:fill="F7F7F8"/></w:rPr><w:t xml:space="preserve">• Implementar: El sistema debe permitir que el mesero ingrese el pedido del cliente y confirmar los ingredientes disponibles antes de transmitir la orden a la cocina de manera rápida y precisa. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:shd w:val="clear" w:color="auto" w:fill="F7F7F8"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:shd w:val="clear" w:color="auto" w:fill="F7F7F8"/></w:rPr><w:t xml:space="preserve">• Organizar: El sistema debe asignar un número y una hora a cada orden, mostrar el estado y tiempo estimado de entrega en una pantalla visible para el mesero.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:shd w:val="clear" w:color="auto" w:fill="F7F7F8"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:shd w:val="clear" w:color="auto" w:fill="F7F7F8"/></w:rPr><w:t xml:space="preserve"> • Facilitar: El sistema facilita la toma de pedidos por parte del mesero y su envío a la cocina de manera instantánea.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml)
